$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") stores values as text (e.g. "30.508.34", using dots as
# thousands separators), not real numbers. Force text formatting on each
# contiguous block of price cells being updated before writing the new values,
# so Excel does not silently reinterpret them as numeric/date values.

$priceBlock = $ws.Range("D2:D4")
$priceBlock.NumberFormat = "@"
$ws.Range("D2").Value = "30.508.34"
$ws.Range("D3").Value = "1.909.85"
$ws.Range("D4").Value = "1.001"
$priceBlock.Style = "Normal"

$priceBlock = $ws.Range("D6:D10")
$priceBlock.NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D7").Value = "0.4832"
$ws.Range("D8").Value = "0.2890"
$ws.Range("D9").Value = "0.06882"
$ws.Range("D10").Value = "111.13"
$priceBlock.Style = "Normal"

$priceBlock = $ws.Range("D12:D22")
$priceBlock.NumberFormat = "@"
$ws.Range("D12").Value = "1.915.40"
$ws.Range("D13").Value = "0.07567"
$ws.Range("D14").Value = "5.381"
$ws.Range("D15").Value = "0.6694"
$ws.Range("D16").Value = "291.87"
$ws.Range("D17").Value = "30.517.57"
$ws.Range("D18").Value = "12.99"
$ws.Range("D19").Value = "1.002"
$ws.Range("D20").Value = "0.000007598"
$ws.Range("D21").Value = "2.166.58"
$ws.Range("D22").Value = "5.503"
$priceBlock.Style = "Normal"

$priceBlock = $ws.Range("D24:D34")
$priceBlock.NumberFormat = "@"
$ws.Range("D24").Value = "6.375"
$ws.Range("D25").Value = "9.446"
$ws.Range("D26").Value = "164.72"
$ws.Range("D27").Value = "20.23"
$ws.Range("D28").Value = "2.080"
$ws.Range("D29").Value = "0.1069"
$ws.Range("D30").Value = "1.439"
$ws.Range("D31").Value = "4.122"
$ws.Range("D32").Value = "4.038"
$ws.Range("D33").Value = "0.04980"
$ws.Range("D34").Value = "0.7359"
$priceBlock.Style = "Normal"

$priceBlock = $ws.Range("D36:D37")
$priceBlock.NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("D37").Value = "2.713"
$priceBlock.Style = "Normal"

$priceBlock = $ws.Range("D39:D44")
$priceBlock.NumberFormat = "@"
$ws.Range("D39").Value = "2.675"
$ws.Range("D40").Value = "2.011"
$ws.Range("D41").Value = "109.52"
$ws.Range("D42").Value = "0.4422"
$ws.Range("D43").Value = "0.8613"
$ws.Range("D44").Value = "5.794"
$priceBlock.Style = "Normal"

$priceBlock = $ws.Range("D46:D51")
$priceBlock.NumberFormat = "@"
$ws.Range("D46").Value = "68.85"
$ws.Range("D47").Value = "7.188"
$ws.Range("D48").Value = "9.215"
$ws.Range("D49").Value = "47.99"
$ws.Range("D50").Value = "0.1225"
$ws.Range("D51").Value = "0.2517"
$priceBlock.Style = "Normal"

# Column E ("Volume(1h)") values are plain percentage-change text, already safe
# to assign directly as strings.
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("E3").Value = "  -0.55%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("E7").Value = "  +1.99%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("E9").Value = "  +1.05%  "
$ws.Range("E10").Value = "  +5.90%  "
$ws.Range("E11").Value = "  +4.92%  "
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("E13").Value = "  -1.72%  "
$ws.Range("E14").Value = "  +2.52%  "
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("E18").Value = "  +0.58%  "
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("E27").Value = "  -3.67%  "
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  +2.91%  "
$ws.Range("E31").Value = "  -1.03%  "
$ws.Range("E32").Value = "  -0.45%  "
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("E38").Value = "  -1.85%  "
$ws.Range("E39").Value = "  -0.50%  "
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("E41").Value = "  -1.22%  "
$ws.Range("E42").Value = "  +0.73%  "
$ws.Range("E43").Value = "  -1.68%  "
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("E46").Value = "  +2.58%  "
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("E48").Value = "  -0.81%  "
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("E51").Value = "  +2.64%  "
